$d = $word.ActiveDocument

# --- Locate the exact range occupied by the old text "Colocación de audífono" ---
$old = "Colocación de audífono"
$new = "Control periódico de la audición"

$search = $d.Content
$search.Find.ClearFormatting()
$search.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if (-not $search.Find.Found) {
    throw "Could not locate the text to replace."
}

$start = $search.Start
$end = $search.End

# --- Replace the whole phrase, then retype the first character separately so
#     that it ends up as its own run (mirrors how the phrase was edited by hand). ---
$range = $d.Range($start, $end)
$range.Text = $new

$splitPos = $start + 1

# A transient bookmark forces the run boundary to persist between "C" and the
# rest of the new sentence, matching the structure of the authored edit.
$d.Bookmarks.Add("_tmpSplit", $d.Range($splitPos, $splitPos)) | Out-Null
$d.Bookmarks("_tmpSplit").Delete()

# --- Move the "_GoBack" bookmark to sit right after the newly typed text,
#     which is where Word leaves it after the last edit. ---
$newEnd = $start + $new.Length
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($newEnd, $newEnd)) | Out-Null
